$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2584.5454
$ws.Range("J70").Value = 2929
$ws.Range("L70").Value = 8787
$ws.Range("N70").Value = -9327
$ws.Range("H73").Value = 2584.5454
$ws.Range("J73").Value = 2929
$ws.Range("L73").Value = 8787
$ws.Range("N73").Value = -10659
$ws.Range("H74").Value = 5546.467
$ws.Range("I74").Value = 5546.467
$ws.Range("K74").Value = 5546.467
$ws.Range("M74").Value = -4610.467
$ws.Range("H77").Value = 5546.467
$ws.Range("I77").Value = 5546.467
$ws.Range("K77").Value = 27732.335
$ws.Range("M77").Value = -23052.335
$ws.Range("H80").Value = 1175.8182
$ws.Range("I80").Value = 1744.2222
$ws.Range("J80").Value = 962.6667
$ws.Range("K80").Value = 5232.6666
$ws.Range("L80").Value = 2888.0001
$ws.Range("M80").Value = -4234.6666
$ws.Range("N80").Value = -4884.0001
$ws.Range("H83").Value = 1175.8182
$ws.Range("I83").Value = 1744.2222
$ws.Range("J83").Value = 962.6667
$ws.Range("K83").Value = 15697.9998
$ws.Range("L83").Value = 8664.0003
$ws.Range("M83").Value = -10705.9998
$ws.Range("N83").Value = -18648.0003
$ws.Range("H112").Value = 2582.92
$ws.Range("J112").Value = 2758.647
$ws.Range("L112").Value = 8275.940999999999
$ws.Range("N112").Value = -10491.941
$ws.Range("H113").Value = 2984.7144
$ws.Range("I113").Value = 2698.75
$ws.Range("K113").Value = 2698.75
$ws.Range("M113").Value = 555.25
$ws.Range("H137").Value = 43489064
$ws.Range("I137").Value = 76925210
$ws.Range("J137").Value = 22080
$ws.Range("K137").Value = 230775630
$ws.Range("L137").Value = 66240
$ws.Range("M137").Value = -230773080
$ws.Range("N137").Value = -71340
$ws.Range("H138").Value = 2321.9822
$ws.Range("J138").Value = 2596.8914
$ws.Range("L138").Value = 7790.674199999999
$ws.Range("N138").Value = -18070.6742
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 158460.3
$ws.Range("I32").Value = 158460.3
$ws.Range("K32").Value = 158460.3
$ws.Range("M32").Value = -158173.3
$ws.Range("H45").Value = 2009.5
$ws.Range("I45").Value = 2133
$ws.Range("J45").Value = 1824.25
$ws.Range("K45").Value = 2133
$ws.Range("L45").Value = 1824.25
$ws.Range("M45").Value = -1756
$ws.Range("N45").Value = -2578.25
$ws.Range("H74").Value = 2791414
$ws.Range("I74").Value = 6175362
$ws.Range("J74").Value = 22729
$ws.Range("K74").Value = 6175362
$ws.Range("L74").Value = 22729
$ws.Range("M74").Value = -6174488
$ws.Range("N74").Value = -24477
$ws.Range("H77").Value = 2791414
$ws.Range("I77").Value = 6175362
$ws.Range("J77").Value = 22729
$ws.Range("K77").Value = 30876810
$ws.Range("L77").Value = 113645
$ws.Range("M77").Value = -30872442
$ws.Range("N77").Value = -122381
$ws.Range("H122").Value = 1393.2222
$ws.Range("I122").Value = 1037.2727
$ws.Range("J122").Value = 1952.5714
$ws.Range("K122").Value = 3111.8181
$ws.Range("L122").Value = 5857.7142
$ws.Range("M122").Value = -661.8181
$ws.Range("N122").Value = -10757.7142
$ws.Range("H132").Value = 1251891.5
$ws.Range("I132").Value = 1390512.9
$ws.Range("K132").Value = 4171538.7
$ws.Range("M132").Value = -4169008.7
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 59182.445
$ws.Range("I20").Value = 75621.92999999999
$ws.Range("J20").Value = 1644.25
$ws.Range("K20").Value = 75621.92999999999
$ws.Range("L20").Value = 1644.25
$ws.Range("M20").Value = -75374.92999999999
$ws.Range("N20").Value = -2138.25
$ws.Range("H99").Value = 35333.332
$ws.Range("I99").Value = 100000
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 100000
$ws.Range("L99").Value = 3000
$ws.Range("M99").Value = -98502
$ws.Range("N99").Value = -5996
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 1000
$ws.Range("M107").Value = 920
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 7000
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 7000
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 7000
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -7340
$ws.Range("H22").Value = 2799.4
$ws.Range("H31").Value = 1686390.8
$ws.Range("I31").Value = 3089359.8
$ws.Range("J31").Value = 2828
$ws.Range("K31").Value = 3089359.8
$ws.Range("L31").Value = 2828
$ws.Range("M31").Value = -3089064.8
$ws.Range("N31").Value = -3418
$ws.Range("H34").Value = 1686390.8
$ws.Range("I34").Value = 3089359.8
$ws.Range("J34").Value = 2828
$ws.Range("K34").Value = 3089359.8
$ws.Range("L34").Value = 2828
$ws.Range("M34").Value = -3089157.8
$ws.Range("N34").Value = -3232
$ws.Range("H68").Value = 0
$ws.Range("I68").Value = 0
$ws.Range("K68").Value = 0
$ws.Range("M68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("I71").Value = 0
$ws.Range("K71").Value = 0
$ws.Range("M71").ClearContents()
$ws.Range("H86").Value = 35166.332
$ws.Range("I86").Value = 35166.332
$ws.Range("K86").Value = 35166.332
$ws.Range("M86").Value = -34043.332
$ws.Range("H89").Value = 35166.332
$ws.Range("I89").Value = 35166.332
$ws.Range("K89").Value = 175831.66
$ws.Range("M89").Value = -170215.66
$ws.Range("H99").Value = 46299.2
$ws.Range("I99").Value = 74499.664
$ws.Range("J99").Value = 3998.5
$ws.Range("K99").Value = 74499.664
$ws.Range("L99").Value = 3998.5
$ws.Range("M99").Value = -73001.664
$ws.Range("N99").Value = -6994.5
$ws.Range("H122").Value = 11221.35
$ws.Range("I122").Value = 2938.4666
$ws.Range("J122").Value = 36070
$ws.Range("K122").Value = 8815.399800000001
$ws.Range("L122").Value = 108210
$ws.Range("M122").Value = -6365.399800000001
$ws.Range("N122").Value = -113110
$ws.Range("H126").Value = 46299.2
$ws.Range("I126").Value = 74499.664
$ws.Range("J126").Value = 3998.5
$ws.Range("K126").Value = 223498.992
$ws.Range("L126").Value = 11995.5
$ws.Range("M126").Value = -221028.992
$ws.Range("N126").Value = -16935.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 2634256.8
$ws.Range("J121").Value = 3849721.8
$ws.Range("L121").Value = 11549165.4
$ws.Range("N121").Value = -11551785.4
$ws.Range("H136").Value = 7218.1665
$ws.Range("I136").Value = 3369.5715
$ws.Range("J136").Value = 12606.2
$ws.Range("K136").Value = 10108.7145
$ws.Range("L136").Value = 37818.60000000001
$ws.Range("M136").Value = -5008.7145
$ws.Range("N136").Value = -48018.60000000001
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 3494
$ws.Range("I113").Value = 2452.4
$ws.Range("K113").Value = 2452.4
$ws.Range("M113").Value = -282.4000000000001
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 7000
$ws.Range("J11").Value = 7000
$ws.Range("L11").Value = 7000
$ws.Range("N11").Value = -7280
$ws.Range("H43").Value = 206000
$ws.Range("H61").Value = 14064.333
$ws.Range("I61").Value = 13330.25
$ws.Range("K61").Value = 13330.25
$ws.Range("M61").Value = -13128.25
$ws.Range("H109").Value = 50000
$ws.Range("J109").Value = 50000
$ws.Range("L109").Value = 50000
$ws.Range("N109").Value = -52774
$ws.Range("H113").Value = 14064.333
$ws.Range("I113").Value = 13330.25
$ws.Range("K113").Value = 13330.25
$ws.Range("M113").Value = -11160.25
$ws.Range("H122").Value = 4749.75
$ws.Range("I122").Value = 4749.75
$ws.Range("K122").Value = 14249.25
$ws.Range("M122").Value = -11799.25
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H37").Value = 12249.5
$ws.Range("I37").Value = 7999.3335
$ws.Range("J37").Value = 25000
$ws.Range("K37").Value = 7999.3335
$ws.Range("L37").Value = 25000
$ws.Range("M37").Value = -7796.3335
$ws.Range("N37").Value = -25406
$ws.Range("H107").Value = 2179.913
$ws.Range("I107").Value = 1340.25
$ws.Range("J107").Value = 7777.6665
$ws.Range("K107").Value = 4020.75
$ws.Range("L107").Value = 23332.9995
$ws.Range("M107").Value = -2100.75
$ws.Range("N107").Value = -27172.9995
$ws.Range("H113").Value = 902.06665
$ws.Range("I113").Value = 1284.875
$ws.Range("J113").Value = 464.57144
$ws.Range("K113").Value = 3854.625
$ws.Range("L113").Value = 1393.71432
$ws.Range("M113").Value = -1684.625
$ws.Range("N113").Value = -5733.71432
$ws.Range("H122").Value = 67788.88
$ws.Range("I122").Value = 1673
$ws.Range("J122").Value = 563658
$ws.Range("K122").Value = 5019
$ws.Range("L122").Value = 1690974
$ws.Range("M122").Value = -2569
$ws.Range("N122").Value = -1695874
$ws.Range("H132").Value = 3877991.8
$ws.Range("I132").Value = 4506278
$ws.Range("K132").Value = 13518834
$ws.Range("M132").Value = -13516304
